$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6756.727
$ws.Range("I40").Value = 6498.8
$ws.Range("J40").Value = 6971.6665
$ws.Range("K40").Value = 6498.8
$ws.Range("L40").Value = 6971.6665
$ws.Range("M40").Value = -6323.8
$ws.Range("N40").Value = -7321.6665

$ws.Range("H62").Value = 4933470
$ws.Range("I62").Value = 5721469
$ws.Range("J62").Value = 8473.75
$ws.Range("K62").Value = 5721469
$ws.Range("L62").Value = 8473.75
$ws.Range("M62").Value = -5720845
$ws.Range("N62").Value = -9721.75

$ws.Range("H65").Value = 4933470
$ws.Range("I65").Value = 5721469
$ws.Range("J65").Value = 8473.75
$ws.Range("K65").Value = 28607345
$ws.Range("L65").Value = 42368.75
$ws.Range("M65").Value = -28604225
$ws.Range("N65").Value = -48608.75

$ws.Range("H82").Value = 9999
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 9999
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 29997
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -30809

$ws.Range("H85").Value = 9999
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 9999
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 29997
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -32805

$ws.Range("H137").Value = 2179.32
$ws.Range("I137").Value = 1361.1428
$ws.Range("K137").Value = 4083.4284
$ws.Range("M137").Value = -1533.4284

$ws.Range("H138").Value = 4138.2
$ws.Range("I138").Value = 1249
$ws.Range("J138").Value = 5293.88
$ws.Range("K138").Value = 3747
$ws.Range("L138").Value = 15881.64
$ws.Range("M138").Value = 1393
$ws.Range("N138").Value = -26161.64

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H27").Value = 6004
$ws.Range("J27").Value = 6004
$ws.Range("L27").Value = 6004
$ws.Range("N27").Value = -6372

$ws.Range("H45").Value = 1238.9048
$ws.Range("I45").Value = 1204.7646
$ws.Range("J45").Value = 1384
$ws.Range("K45").Value = 1204.7646
$ws.Range("L45").Value = 1384
$ws.Range("M45").Value = -827.7646
$ws.Range("N45").Value = -2138

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 21302.375
$ws.Range("J74").Value = 21302.375
$ws.Range("L74").Value = 21302.375
$ws.Range("N74").Value = -23174.375

$ws.Range("H77").Value = 21302.375
$ws.Range("J77").Value = 21302.375
$ws.Range("L77").Value = 63907.125
$ws.Range("N77").Value = -73267.125

$ws.Range("H99").Value = 6808.02
$ws.Range("I99").Value = 6562.3555
$ws.Range("J99").Value = 9019
$ws.Range("K99").Value = 6562.3555
$ws.Range("L99").Value = 9019
$ws.Range("M99").Value = -5064.3555
$ws.Range("N99").Value = -12015

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1998.25
$ws.Range("J31").Value = 2955
$ws.Range("L31").Value = 2955
$ws.Range("N31").Value = -3545

$ws.Range("H34").Value = 1998.25
$ws.Range("J34").Value = 2955
$ws.Range("L34").Value = 2955
$ws.Range("N34").Value = -3359

$ws.Range("H58").Value = 41671604
$ws.Range("I58").Value = 47622030
$ws.Range("K58").Value = 47622030
$ws.Range("M58").Value = -47621827

$ws.Range("H62").Value = 13874.833
$ws.Range("I62").Value = 8875
$ws.Range("K62").Value = 8875
$ws.Range("M62").Value = -8251

$ws.Range("H65").Value = 13874.833
$ws.Range("I65").Value = 8875
$ws.Range("K65").Value = 44375
$ws.Range("M65").Value = -41255

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H132").Value = 3630.9714
$ws.Range("I132").Value = 2677.5862
$ws.Range("K132").Value = 8032.758600000001
$ws.Range("M132").Value = -5502.758600000001

$ws.Range("H136").Value = 41671604
$ws.Range("I136").Value = 47622030
$ws.Range("K136").Value = 142866090
$ws.Range("M136").Value = -142863540

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 162281
$ws.Range("I5").Value = 915.13043
$ws.Range("J5").Value = 626207.9
$ws.Range("K5").Value = 2745.39129
$ws.Range("L5").Value = 1878623.7
$ws.Range("M5").Value = -2633.39129
$ws.Range("N5").Value = -1878847.7

$ws.Range("H37").Value = 218174.55
$ws.Range("J37").Value = 218174.55
$ws.Range("L37").Value = 654523.6499999999
$ws.Range("N37").Value = -654747.6499999999

$ws.Range("H76").Value = 10197.125
$ws.Range("I76").Value = 8096.7144
$ws.Range("J76").Value = 24900
$ws.Range("K76").Value = 24290.1432
$ws.Range("L76").Value = 74700
$ws.Range("M76").Value = -23907.1432
$ws.Range("N76").Value = -75466

$ws.Range("H79").Value = 10197.125
$ws.Range("I79").Value = 8096.7144
$ws.Range("J79").Value = 24900
$ws.Range("K79").Value = 24290.1432
$ws.Range("L79").Value = 74700
$ws.Range("M79").Value = -22964.1432
$ws.Range("N79").Value = -77352

$ws.Range("H80").Value = 4166
$ws.Range("I80").Value = 1999
$ws.Range("K80").Value = 5997
$ws.Range("M80").Value = -5061

$ws.Range("H83").Value = 4166
$ws.Range("I83").Value = 1999
$ws.Range("K83").Value = 17991
$ws.Range("M83").Value = -13311

$ws.Range("H87").Value = 23302
$ws.Range("I87").Value = 11371.333
$ws.Range("J87").Value = 32250
$ws.Range("K87").Value = 34113.999
$ws.Range("L87").Value = 96750
$ws.Range("M87").Value = -32865.999
$ws.Range("N87").Value = -99246

$ws.Range("H90").Value = 23302
$ws.Range("I90").Value = 11371.333
$ws.Range("J90").Value = 32250
$ws.Range("K90").Value = 102341.997
$ws.Range("L90").Value = 290250
$ws.Range("M90").Value = -96101.997
$ws.Range("N90").Value = -302730

$ws.Range("H111").Value = 6064.5
$ws.Range("I111").Value = 2357.6667
$ws.Range("K111").Value = 7073.000100000001
$ws.Range("M111").Value = -4006.000100000001

$ws.Range("H121").Value = 131286
$ws.Range("J121").Value = 187237.14
$ws.Range("L121").Value = 561711.42
$ws.Range("N121").Value = -564331.42

$ws.Range("H134").Value = 4789.2
$ws.Range("J134").Value = 7033
$ws.Range("L134").Value = 21099
$ws.Range("N134").Value = -31239

$ws.Range("H135").Value = 162281
$ws.Range("I135").Value = 915.13043
$ws.Range("J135").Value = 626207.9
$ws.Range("K135").Value = 8236.173870000001
$ws.Range("L135").Value = 5635871.100000001
$ws.Range("M135").Value = -5701.173870000001
$ws.Range("N135").Value = -5640941.100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 55558790
$ws.Range("I132").Value = 71431060
$ws.Range("K132").Value = 214293180
$ws.Range("M132").Value = -214290650

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5438.65
$ws.Range("J7").Value = 5953.2666
$ws.Range("L7").Value = 5953.2666
$ws.Range("N7").Value = -6177.2666

$ws.Range("H55").Value = 2436.4285
$ws.Range("J55").Value = 4117.091
$ws.Range("L55").Value = 4117.091
$ws.Range("N55").Value = -4463.091

$ws.Range("H61").Value = 6861.0835
$ws.Range("I61").Value = 5283.7144
$ws.Range("J61").Value = 9069.4
$ws.Range("K61").Value = 5283.7144
$ws.Range("L61").Value = 9069.4
$ws.Range("M61").Value = -5081.7144
$ws.Range("N61").Value = -9473.4

$ws.Range("H113").Value = 6861.0835
$ws.Range("I113").Value = 5283.7144
$ws.Range("J113").Value = 9069.4
$ws.Range("K113").Value = 5283.7144
$ws.Range("L113").Value = 9069.4
$ws.Range("M113").Value = -3113.7144
$ws.Range("N113").Value = -13409.4

$ws.Range("H126").Value = 5438.65
$ws.Range("J126").Value = 5953.2666
$ws.Range("L126").Value = 17859.7998
$ws.Range("N126").Value = -22799.7998

$ws.Range("H136").Value = 15626983
$ws.Range("I136").Value = 27778628
$ws.Range("K136").Value = 83335884
$ws.Range("M136").Value = -83333334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 4450
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

$ws.Range("H17").Value = 905
$ws.Range("I17").Value = 905
$ws.Range("K17").Value = 905
$ws.Range("M17").Value = -733

$ws.Range("H107").Value = 367.7
$ws.Range("I107").Value = 367.7
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 367.7
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 816.9000000000001
$ws.Range("N107").ClearContents()

Write-Host "Applied all changes"
